$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Eintragungsdatum" value from 15.12.2023 to 01.01.2024
$ws.Range("B7").Value = "01.01.2024"

# Move the active selection to B4 (as reflected in the saved file)
$ws.Range("B4").Select()
